$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (Arkansas) updates ---
$ws.Cells.Item(9, 3).Value = 24253
$ws.Cells.Item(9, 4).Value = 292
$ws.Cells.Item(9, 5).Value = 5239
$ws.Cells.Item(9, 7).Value = 25.23
$ws.Cells.Item(9, 8).Value = 25.18
$ws.Cells.Item(9, 11).Value = 20762
$ws.Cells.Item(9, 12).Value = 278

# --- Row 10 (California - San Diego) updates ---
$ws.Cells.Item(10, 2).Value = 44018
$ws.Cells.Item(10, 3).Value = 17000
$ws.Cells.Item(10, 5).Value = 615
$ws.Cells.Item(10, 7).Value = 4.55
$ws.Cells.Item(10, 11).Value = 13514

# --- Row 31 (Washington) updates ---
$ws.Cells.Item(31, 2).Value = 44018
$ws.Cells.Item(31, 3).Value = 36985
$ws.Cells.Item(31, 4).Value = 1370
$ws.Cells.Item(31, 5).Value = 1455
$ws.Cells.Item(31, 11).Value = 26515

# --- New row 42 (Iowa) ---
$ws.Cells.Item(42, 1).Value = "Iowa"
$ws.Cells.Item(42, 9).Value = $false
$ws.Cells.Item(42, 10).Value = $false
$ws.Cells.Item(42, 13).Value = 109911
$ws.Cells.Item(42, 14).Value = 3.51
$q = [char]34
$ws.Cells.Item(42, 15).Value = "An error occurred. ... ValueError('Unable to parse $q" + "Reported Deaths In Adair : No Data" + "$q as int')"
